# Generate Report for handback
# Update the Correspond Handoff Datetime (column D) and
# Correspond Handback DateTime (column G) for the "11330032" file rows
# on the zh-cn and de-de report sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-02-16 09:54:07"
$wsZhCn.Range("G2").Value = "2016-02-16 09:55:01"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-02-16 09:54:21"
$wsDeDe.Range("G2").Value = "2016-02-16 09:55:29"
